$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the NCI ID values that referenced the retired oq512c512h / oq512c512ht
# shared strings (job-name column, no longer recorded for these runs)
$ws.Range("B2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()

# Record a Duration value for the run in row 5
$ws.Cells.Item(5, 9).Value = 16

# New job tracking row for the PGA / param file run
$ws.Cells.Item(7, 1).Value = "Domains merged collapsed"
$ws.Cells.Item(7, 2).Value = 5039046
$ws.Cells.Item(7, 3).Value = 0.41666666666666669
$ws.Cells.Item(7, 4).Value = 256
$ws.Cells.Item(7, 5).Value = 512
$ws.Cells.Item(7, 6).Value = 54
$ws.Cells.Item(7, 7).Value = 200
$ws.Cells.Item(7, 8).Value = 4.5
$ws.Cells.Item(7, 9).Value = 1.7
$ws.Cells.Item(7, 10).Value = "Complete"
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = "PGA"
$ws.Cells.Item(7, 13).Value = 15

# Match formatting used by the other data rows: Walltime column is a time
# value, CPUs column uses the Menlo font when called out, Duration column
# is shown to one decimal place
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("I6").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection below the newly added row
$ws.Range("A8").Select()
